$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (tab name within workbook, distinct from the sheet's internal name)
$ws.Name = "tot_co2_hist"

$values = @{
    "G2" = 4
    "H2" = 11
    "K2" = 39
    "L2" = 126
    "M2" = 363
    "N2" = 1605
    "O2" = 6707
    "P2" = 24040
    "Q2" = 83721
    "R2" = 197920
    "S2" = 329782
    "T2" = 670822
    "U2" = 566388
    "V2" = 681906
    "W2" = 502067
    "X2" = 580795
    "Y2" = 558301
    "Z2" = 565934
    "AA2" = 520330
    "AB2" = 620310
    "AC2" = 602223
    "AD2" = 726874
    "AE2" = 780596
    "AF2" = 875704
    "AG2" = 1035368
    "AH2" = 1254865
    "AI2" = 1394379
    "AJ2" = 1773189
    "AK2" = 1718216
    "AL2" = 1965757
    "AM2" = 2324388
    "AN2" = 2536766
    "AO2" = 2609083
    "AP2" = 2736295
    "AQ2" = 2728257
    "AR2" = 3026526
    "AS2" = 2941449
    "AT2" = 2829510
    "AU2" = 3253643
    "AV2" = 3201100
    "AW2" = 2908431
    "AX2" = 2755344
    "AY2" = 2627556
    "AZ2" = 2789685
    "BA2" = 2597067
    "BB2" = 2463421
    "BC2" = 2286219
    "BD2" = 2209485
    "BE2" = 2152019
    "BF2" = 2023890
    "BG2" = 1987335
    "BH2" = 1973836
    "BI2" = 2277615
    "BJ2" = 2173150
    "BK2" = 2410007
    "BL2" = 2119212
    "BM2" = 2088742
    "BN2" = 2134846
    "BO2" = 2086935
    "BP2" = 1886528
    "BQ2" = 1994147
    "BR2" = 2298845
    "BS2" = 2373184
    "BT2" = 2252170
    "BU2" = 2347090
    "BV2" = 2364229
    "BW2" = 2271242
    "BX2" = 2710381
    "BY2" = 2590410
    "BZ2" = 2454336
    "CA2" = 2505217
    "CB2" = 2408268
    "CC2" = 2410351
    "CD2" = 2235857
    "CE2" = 2215073
    "CF2" = 2168207
    "CG2" = 2098984
    "CH2" = 1949830
    "CI2" = 2443591
    "CJ2" = 2115200
    "CK2" = 1938512
    "CL2" = 1742123
    "CM2" = 1577263
    "CN2" = 1525663
    "CO2" = 1405365
    "CP2" = 1328098
    "CQ2" = 1200262
    "CR2" = 1100109
    "CS2" = 998085
    "CT2" = 890817
    "CU2" = 847545
    "CV2" = 726816
    "CW2" = 615967
    "CX2" = 525429
    "CY2" = 461469
    "CZ2" = 376736
    "DA2" = 334423
    "DB2" = 277732
    "DC2" = 226476
    "DD2" = 187602
    "DE2" = 155765
    "DF2" = 123013
    "DG2" = 101457
    "DH2" = 81331
    "DI2" = 64180
    "DJ2" = 51568
    "DK2" = 41880
    "DL2" = 34422
    "DM2" = 28555
    "DN2" = 23798
    "DO2" = 23003
    "DP2" = 17646
    "DQ2" = 17645
    "DR2" = 15936
    "DS2" = 14117
    "DT2" = 12156
    "DU2" = 9127
    "DV2" = 7734
    "DW2" = 5819
    "DX2" = 4130
    "DY2" = 2965
    "DZ2" = 1991
    "EA2" = 1477
    "EB2" = 1137
    "EC2" = 776
    "ED2" = 646
    "EE2" = 477
    "EF2" = 348
    "EG2" = 253
    "EH2" = 200
    "EI2" = 141
    "EJ2" = 83
    "EK2" = 56
    "EL2" = 51
    "EN2" = 21
    "EO2" = 23
    "EQ2" = 4
}

foreach ($cellRef in $values.Keys) {
    $ws.Range($cellRef).Value = $values[$cellRef]
}
